$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Enter new text values in the precise order the original author typed
#    them, so that the shared-strings table is rebuilt with the same slot
#    ordering as the authoritative edit (the freed slot from the old D4
#    text is reused by the first brand new string "FREND4").
# ---------------------------------------------------------------------------
$ws.Range("A12").Value2 = "FREND4"
$ws.Range("C12").Value2 = "Add a dropdown of the titles or the primary key field when doing an update. This can help in preventing creation of a new record."
$ws.Range("B12").Value2 = "Updates can create a new record if not properly entered."
$ws.Range("D4").Value2  = "Author position is now removed hence bug fix is not required."
$ws.Range("D6").Value2  = "Report is now generated in groups hence bug fix is not required."

# ---------------------------------------------------------------------------
# 2. Apply formatting. Use copy/paste-special of formats from cells that
#    already carry the desired style (the wrap-text style shared by the
#    other long-text cells in columns B/C/D) so the engine reuses the
#    existing cellXfs entry instead of synthesizing a brand-new one.
# ---------------------------------------------------------------------------
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# D6 should look like the other wrap-text "Potential Fix"/"Status" cells
$ws.Range("C4").Copy() | Out-Null
$ws.Range("D6").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# C12 (Potential Fix column) gets the same wrap-text style
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C12").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# A12 / B12 stay in the plain default style (matches other Code / Bug
# Description cells such as A5 / B5), nothing further required.

# ---------------------------------------------------------------------------
# 3. Row height for the newly added row, matching the wrapped-text row
#    height used by similarly sized entries (row 4 is also 57.6pt).
# ---------------------------------------------------------------------------
$ws.Rows(12).RowHeight = 57.6

# ---------------------------------------------------------------------------
# 4. Update the active selection to the last edited cell, like Excel would
#    leave it after the final edit.
# ---------------------------------------------------------------------------
$ws.Range("D12").Select() | Out-Null
